# Add the "mini PEKKA" team row (new row 4), duplicating the existing
# row 3 values, matching the diff which appends this row and expands
# the sheet dimension to A1:F4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 'mini PEKKA'
$ws.Range("B4").Value = 'Federico Leonardi | Rita Levi’s'
$ws.Range("C4").Value = 'Davide  Rosa` | Hellas Lazio'
$ws.Range("D4").Value = 'Lorenzo Casari | Nazzzionale ferrovieri'
$ws.Range("E4").Value = 'Mattia Festi | SHARK ATTACK'
$ws.Range("F4").Value = 'Matteo Gatti | demobusters'
